# Guion 04 grado 03 - seguimiento update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grado: value changes from 6 to 3
$ws.Range("C2").Value = 3

# Row 9 (item 4) gets dated / reviewed:
#   B9: date 2015-03-12 (serial 42075) with a heavier left border (no right border)
#   C9, D9, E9: same date, keep existing date-style formatting
#   F9: observation text (new shared string)
$ws.Range("C9").Value = 42075
$ws.Range("D9").Value = 42075
$ws.Range("E9").Value = 42075
$ws.Range("F9").Value = "En revisión por parte de editor"

# B9 needs its own look: keep the date format/fill, drop the right border
$ws.Range("B9").Value = 42075
$ws.Range("B9").Borders.Item(10).LineStyle = -4142

# Row 9 grows taller to fit the wrapped observation text
$ws.Rows.Item(9).RowHeight = 30.75

# Selection moves to B9 (and the view no longer needs to be scrolled to A4)
$null = $ws.Range("B9").Select()

Write-Host "done"
